$d = $word.ActiveDocument

# --- Part 1: append the new second paragraph -----------------------------
# Do this first, while the first paragraph's run still has no explicit
# character formatting, so the freshly inserted paragraph does not inherit
# any stray run properties from part 2 below.
$para1 = $d.Paragraphs(1)
$para1.Range.InsertParagraphAfter()
$para2 = $d.Paragraphs(2)
$para2.Range.Text = "je le change maintenant et je regarde ce qui se passe"

# --- Part 2: split "word" out of the first paragraph's run ---------------
# Before:
#   "je crée un fichier word pour étudier son suivi avec git"   (one run)
# After:
#   "je crée un fichier " + "word" + " pour étudier son suivi avec git"
# as three separate runs (the middle one is where Word's proofing pass
# flags the English word "word" amid French text). Locate the exact
# character span of "word" inside the paragraph and force Word to
# materialize a run boundary there.
$para1 = $d.Paragraphs(1)
$fullText = $para1.Range.Text
$target = "word"
$startOffset = $fullText.IndexOf($target)
$paraStart = $para1.Range.Start
$wordStart = $paraStart + $startOffset
$wordEnd = $wordStart + $target.Length

$wordRange = $d.Range($wordStart, $wordEnd)
$wordRange.Bold = $true
$wordRange.Bold = $false
